$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 4 currently reads "中雨，今天是农历五月初四，明天又是端午节了"
# with a paragraph-mark <w:rFonts w:hint="default"/> and carries the
# "_GoBack" bookmark. The edit:
#   1. Flips that paragraph's mark hint to "eastAsia" and strips the
#      bookmark (it is going to move onto the new last paragraph).
#   2. Inserts a new paragraph "2022年6月3日星期五" right after it.
#   3. Inserts a further new paragraph "中雨，今天是农历五月初五，中国
#      传统端午节" after that, which now owns the "_GoBack" bookmark.
# The trailing empty paragraph stays put at the end of the story.
# ---------------------------------------------------------------------------

$wdNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Step 1: rewrite paragraph 4 (drop the bookmark, fix the hint) --------
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range.Duplicate
$p4Xml = '<w:p ' + $wdNS + '><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>&#20013;&#38632;&#65292;&#20170;&#22825;&#26159;&#20892;&#21382;&#20116;&#26376;&#21021;&#22235;&#65292;&#26126;&#22825;&#21448;&#26159;&#31471;&#21320;&#33410;&#20102;</w:t></w:r></w:p>'
$r4.InsertXML($p4Xml)

# --- Step 2: add "2022年6月3日星期五" right after paragraph 4 ------------
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range.Duplicate
$r5.Collapse(1)
$p5Xml = '<w:p ' + $wdNS + '><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>2022&#24180;6&#26376;3&#26085;&#26143;&#26399;&#20116;</w:t></w:r></w:p>'
$r5.InsertXML($p5Xml)

# --- Step 3: add "中雨，今天是农历五月初五，中国传统端午节" + bookmark ---
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()

$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range.Duplicate
$r6.Collapse(1)
$p6Xml = '<w:p ' + $wdNS + '><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>&#20013;&#38632;&#65292;&#20170;&#22825;&#26159;&#20892;&#21382;&#20116;&#26376;&#21021;&#20116;&#65292;&#20013;&#22269;&#20256;&#32479;&#31471;&#21320;&#33410;</w:t></w:r></w:p>'
$r6.InsertXML($p6Xml)

$p6 = $d.Paragraphs.Item(6)
$r6b = $p6.Range.Duplicate
$r6b.Collapse(0)
$r6b.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $r6b)

Write-Host "Final paragraph count: " $d.Paragraphs.Count
